$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-17 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-18 Saturday", 2) | Out-Null
$d.Content.Find.Execute("36+17=53", $true, $false, $false, $false, $false, $true, 1, $false, "26-9=17", 2) | Out-Null
$d.Content.Find.Execute("94-94=0", $true, $false, $false, $false, $false, $true, 1, $false, "58-31=27", 2) | Out-Null
$d.Content.Find.Execute("56-27=29", $true, $false, $false, $false, $false, $true, 1, $false, "40-13=27", 2) | Out-Null
$d.Content.Find.Execute("55-21=34", $true, $false, $false, $false, $false, $true, 1, $false, "61-31=30", 2) | Out-Null
$d.Content.Find.Execute("32+18=50", $true, $false, $false, $false, $false, $true, 1, $false, "64+27=91", 2) | Out-Null
$d.Content.Find.Execute("82+4=86", $true, $false, $false, $false, $false, $true, 1, $false, "48+30=78", 2) | Out-Null
$d.Content.Find.Execute("52+2=54", $true, $false, $false, $false, $false, $true, 1, $false, "79-5=74", 2) | Out-Null
$d.Content.Find.Execute("66+29=95", $true, $false, $false, $false, $false, $true, 1, $false, "72-26=46", 2) | Out-Null
$d.Content.Find.Execute("60+18=78", $true, $false, $false, $false, $false, $true, 1, $false, "17+24=41", 2) | Out-Null
$d.Content.Find.Execute("52-25=27", $true, $false, $false, $false, $false, $true, 1, $false, "91-32=59", 2) | Out-Null
$d.Content.Find.Execute("86-29=57", $true, $false, $false, $false, $false, $true, 1, $false, "59+3=62", 2) | Out-Null
$d.Content.Find.Execute("70-38=32", $true, $false, $false, $false, $false, $true, 1, $false, "51-51=0", 2) | Out-Null
$d.Content.Find.Execute("11-0=11", $true, $false, $false, $false, $false, $true, 1, $false, "31-24=7", 2) | Out-Null
$d.Content.Find.Execute("74-48=26", $true, $false, $false, $false, $false, $true, 1, $false, "8+75=83", 2) | Out-Null
$d.Content.Find.Execute("7+31=38", $true, $false, $false, $false, $false, $true, 1, $false, "2+74=76", 2) | Out-Null
$d.Content.Find.Execute("15+46=61", $true, $false, $false, $false, $false, $true, 1, $false, "6+54=60", 2) | Out-Null
$d.Content.Find.Execute("32-13=19", $true, $false, $false, $false, $false, $true, 1, $false, "18+46=64", 2) | Out-Null
$d.Content.Find.Execute("57-31=26", $true, $false, $false, $false, $false, $true, 1, $false, "98-47=51", 2) | Out-Null
$d.Content.Find.Execute("5+36=41", $true, $false, $false, $false, $false, $true, 1, $false, "58-29=29", 2) | Out-Null
$d.Content.Find.Execute("2+82=84", $true, $false, $false, $false, $false, $true, 1, $false, "72-19=53", 2) | Out-Null
$d.Content.Find.Execute("7-2=5", $true, $false, $false, $false, $false, $true, 1, $false, "36+2=38", 2) | Out-Null
$d.Content.Find.Execute("23+74=97", $true, $false, $false, $false, $false, $true, 1, $false, "22+19=41", 2) | Out-Null
$d.Content.Find.Execute("1+78=79", $true, $false, $false, $false, $false, $true, 1, $false, "16-7=9", 2) | Out-Null
$d.Content.Find.Execute("25+73=98", $true, $false, $false, $false, $false, $true, 1, $false, "6+20=26", 2) | Out-Null
$d.Content.Find.Execute("27-20=7", $true, $false, $false, $false, $false, $true, 1, $false, "69+8=77", 2) | Out-Null
$d.Content.Find.Execute("38-35=3", $true, $false, $false, $false, $false, $true, 1, $false, "64-15=49", 2) | Out-Null
$d.Content.Find.Execute("84-73=11", $true, $false, $false, $false, $false, $true, 1, $false, "86-19=67", 2) | Out-Null
$d.Content.Find.Execute("59+6=65", $true, $false, $false, $false, $false, $true, 1, $false, "51-20=31", 2) | Out-Null
$d.Content.Find.Execute("23-18=5", $true, $false, $false, $false, $false, $true, 1, $false, "65-9=56", 2) | Out-Null
$d.Content.Find.Execute("63+34=97", $true, $false, $false, $false, $false, $true, 1, $false, "84-41=43", 2) | Out-Null
$d.Content.Find.Execute("78-42=36", $true, $false, $false, $false, $false, $true, 1, $false, "47+38=85", 2) | Out-Null
$d.Content.Find.Execute("92-80=12", $true, $false, $false, $false, $false, $true, 1, $false, "23+1=24", 2) | Out-Null
$d.Content.Find.Execute("76+14=90", $true, $false, $false, $false, $false, $true, 1, $false, "26+46=72", 2) | Out-Null
$d.Content.Find.Execute("21+6=27", $true, $false, $false, $false, $false, $true, 1, $false, "49+30=79", 2) | Out-Null
$d.Content.Find.Execute("32+51=83", $true, $false, $false, $false, $false, $true, 1, $false, "43+4=47", 2) | Out-Null
$d.Content.Find.Execute("20+56=76", $true, $false, $false, $false, $false, $true, 1, $false, "35+31=66", 2) | Out-Null
$d.Content.Find.Execute("44+24=68", $true, $false, $false, $false, $false, $true, 1, $false, "3+84=87", 2) | Out-Null
$d.Content.Find.Execute("48+6=54", $true, $false, $false, $false, $false, $true, 1, $false, "58-32=26", 2) | Out-Null
$d.Content.Find.Execute("79-67=12", $true, $false, $false, $false, $false, $true, 1, $false, "24+7=31", 2) | Out-Null
$d.Content.Find.Execute("32+6=38", $true, $false, $false, $false, $false, $true, 1, $false, "81-32=49", 2) | Out-Null
$d.Content.Find.Execute("36-15=21", $true, $false, $false, $false, $false, $true, 1, $false, "95-4=91", 2) | Out-Null
$d.Content.Find.Execute("67+29=96", $true, $false, $false, $false, $false, $true, 1, $false, "14+82=96", 2) | Out-Null
$d.Content.Find.Execute("43+8=51", $true, $false, $false, $false, $false, $true, 1, $false, "27-18=9", 2) | Out-Null
$d.Content.Find.Execute("79+16=95", $true, $false, $false, $false, $false, $true, 1, $false, "97-74=23", 2) | Out-Null
$d.Content.Find.Execute("62+15=77", $true, $false, $false, $false, $false, $true, 1, $false, "47+34=81", 2) | Out-Null
$d.Content.Find.Execute("69-47=22", $true, $false, $false, $false, $false, $true, 1, $false, "89-75=14", 2) | Out-Null
$d.Content.Find.Execute("68-60=8", $true, $false, $false, $false, $false, $true, 1, $false, "45-31=14", 2) | Out-Null
$d.Content.Find.Execute("64-12=52", $true, $false, $false, $false, $false, $true, 1, $false, "94-48=46", 2) | Out-Null
$d.Content.Find.Execute("5+63=68", $true, $false, $false, $false, $false, $true, 1, $false, "61+33=94", 2) | Out-Null
$d.Content.Find.Execute("86-15=71", $true, $false, $false, $false, $false, $true, 1, $false, "73-37=36", 2) | Out-Null
$d.Content.Find.Execute("21-17=4", $true, $false, $false, $false, $false, $true, 1, $false, "10+38=48", 2) | Out-Null
$d.Content.Find.Execute("13+66=79", $true, $false, $false, $false, $false, $true, 1, $false, "15+47=62", 2) | Out-Null
$d.Content.Find.Execute("53+7=60", $true, $false, $false, $false, $false, $true, 1, $false, "24+19=43", 2) | Out-Null
$d.Content.Find.Execute("39+42=81", $true, $false, $false, $false, $false, $true, 1, $false, "8+33=41", 2) | Out-Null
$d.Content.Find.Execute("9+25=34", $true, $false, $false, $false, $false, $true, 1, $false, "66+14=80", 2) | Out-Null
$d.Content.Find.Execute("78-41=37", $true, $false, $false, $false, $false, $true, 1, $false, "88-42=46", 2) | Out-Null
$d.Content.Find.Execute("54+19=73", $true, $false, $false, $false, $false, $true, 1, $false, "33+64=97", 2) | Out-Null
$d.Content.Find.Execute("99-31=68", $true, $false, $false, $false, $false, $true, 1, $false, "94-59=35", 2) | Out-Null
$d.Content.Find.Execute("12+18=30", $true, $false, $false, $false, $false, $true, 1, $false, "50+26=76", 2) | Out-Null
$d.Content.Find.Execute("10-10=0", $true, $false, $false, $false, $false, $true, 1, $false, "14+72=86", 2) | Out-Null
$d.Content.Find.Execute("58+21=79", $true, $false, $false, $false, $false, $true, 1, $false, "67+32=99", 2) | Out-Null
$d.Content.Find.Execute("22+11=33", $true, $false, $false, $false, $false, $true, 1, $false, "28-8=20", 2) | Out-Null
$d.Content.Find.Execute("42-1=41", $true, $false, $false, $false, $false, $true, 1, $false, "3+38=41", 2) | Out-Null
$d.Content.Find.Execute("77-71=6", $true, $false, $false, $false, $false, $true, 1, $false, "59-48=11", 2) | Out-Null
$d.Content.Find.Execute("63+26=89", $true, $false, $false, $false, $false, $true, 1, $false, "21+73=94", 2) | Out-Null
$d.Content.Find.Execute("95-69=26", $true, $false, $false, $false, $false, $true, 1, $false, "3+10=13", 2) | Out-Null
$d.Content.Find.Execute("83-16=67", $true, $false, $false, $false, $false, $true, 1, $false, "18+14=32", 2) | Out-Null
$d.Content.Find.Execute("51-16=35", $true, $false, $false, $false, $false, $true, 1, $false, "84-33=51", 2) | Out-Null
$d.Content.Find.Execute("68+29=97", $true, $false, $false, $false, $false, $true, 1, $false, "12+33=45", 2) | Out-Null
$d.Content.Find.Execute("52+18=70", $true, $false, $false, $false, $false, $true, 1, $false, "82-70=12", 2) | Out-Null
$d.Content.Find.Execute("8-7=1", $true, $false, $false, $false, $false, $true, 1, $false, "74+22=96", 2) | Out-Null
$d.Content.Find.Execute("6-4=2", $true, $false, $false, $false, $false, $true, 1, $false, "74-36=38", 2) | Out-Null
$d.Content.Find.Execute("22+46=68", $true, $false, $false, $false, $false, $true, 1, $false, "29+29=58", 2) | Out-Null
$d.Content.Find.Execute("72+19=91", $true, $false, $false, $false, $false, $true, 1, $false, "48+35=83", 2) | Out-Null
$d.Content.Find.Execute("69-60=9", $true, $false, $false, $false, $false, $true, 1, $false, "27-22=5", 2) | Out-Null
$d.Content.Find.Execute("86+11=97", $true, $false, $false, $false, $false, $true, 1, $false, "20-14=6", 2) | Out-Null
$d.Content.Find.Execute("52-49=3", $true, $false, $false, $false, $false, $true, 1, $false, "83-45=38", 2) | Out-Null
$d.Content.Find.Execute("26-21=5", $true, $false, $false, $false, $false, $true, 1, $false, "68-26=42", 2) | Out-Null
$d.Content.Find.Execute("66+6=72", $true, $false, $false, $false, $false, $true, 1, $false, "21+26=47", 2) | Out-Null
$d.Content.Find.Execute("59+25=84", $true, $false, $false, $false, $false, $true, 1, $false, "71+28=99", 2) | Out-Null
$d.Content.Find.Execute("81-66=15", $true, $false, $false, $false, $false, $true, 1, $false, "8+82=90", 2) | Out-Null
$d.Content.Find.Execute("50-25=25", $true, $false, $false, $false, $false, $true, 1, $false, "51+20=71", 2) | Out-Null
$d.Content.Find.Execute("27+49=76", $true, $false, $false, $false, $false, $true, 1, $false, "16+0=16", 2) | Out-Null
$d.Content.Find.Execute("57-7=50", $true, $false, $false, $false, $false, $true, 1, $false, "30+37=67", 2) | Out-Null
$d.Content.Find.Execute("99-10=89", $true, $false, $false, $false, $false, $true, 1, $false, "84-33=51", 2) | Out-Null
$d.Content.Find.Execute("66-42=24", $true, $false, $false, $false, $false, $true, 1, $false, "14+82=96", 2) | Out-Null
$d.Content.Find.Execute("7+34=41", $true, $false, $false, $false, $false, $true, 1, $false, "42+2=44", 2) | Out-Null
$d.Content.Find.Execute("52-11=41", $true, $false, $false, $false, $false, $true, 1, $false, "9+4=13", 2) | Out-Null
$d.Content.Find.Execute("81-30=51", $true, $false, $false, $false, $false, $true, 1, $false, "17+40=57", 2) | Out-Null
$d.Content.Find.Execute("14-0=14", $true, $false, $false, $false, $false, $true, 1, $false, "54+37=91", 2) | Out-Null
$d.Content.Find.Execute("61-36=25", $true, $false, $false, $false, $false, $true, 1, $false, "96-36=60", 2) | Out-Null
$d.Content.Find.Execute("34+25=59", $true, $false, $false, $false, $false, $true, 1, $false, "99-1=98", 2) | Out-Null
$d.Content.Find.Execute("36+57=93", $true, $false, $false, $false, $false, $true, 1, $false, "34+42=76", 2) | Out-Null
$d.Content.Find.Execute("51+24=75", $true, $false, $false, $false, $false, $true, 1, $false, "47+49=96", 2) | Out-Null
$d.Content.Find.Execute("6+39=45", $true, $false, $false, $false, $false, $true, 1, $false, "81-73=8", 2) | Out-Null
$d.Content.Find.Execute("21+8=29", $true, $false, $false, $false, $false, $true, 1, $false, "90-33=57", 2) | Out-Null
$d.Content.Find.Execute("94-2=92", $true, $false, $false, $false, $false, $true, 1, $false, "75-47=28", 2) | Out-Null
$d.Content.Find.Execute("62+1=63", $true, $false, $false, $false, $false, $true, 1, $false, "79-43=36", 2) | Out-Null
$d.Content.Find.Execute("86-44=42", $true, $false, $false, $false, $false, $true, 1, $false, "88-72=16", 2) | Out-Null
$d.Content.Find.Execute("99-90=9", $true, $false, $false, $false, $false, $true, 1, $false, "56-32=24", 2) | Out-Null
